$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C16").Value = "1050952506"
$ws.Range("D16").Value = "YURIS ROCIO PUELLO OYOLA"
$ws.Range("E16").Value = "1712"
$ws.Range("F16").Value = 29509
$ws.Range("G16").Value = 737717
$ws.Range("C17").Value = "92227074"
$ws.Range("D17").Value = "RODRIGO ANTONIO VELASQUEZ MORELO"
$ws.Range("E17").Value = "1712"
$ws.Range("F17").Value = 29509
$ws.Range("G17").Value = 737717
$ws.Range("C18").Value = "1050963935"
$ws.Range("D18").Value = "ANGELA INES CAUSIL MARTINEZ"
$ws.Range("E18").Value = "1712"
$ws.Range("F18").Value = 35730
$ws.Range("G18").Value = 893263
$ws.Range("C19").Value = "1050952506"
$ws.Range("D19").Value = "YURIS ROCIO PUELLO OYOLA"
$ws.Range("E19").Value = "1801"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 737717
$ws.Range("C20").Value = "92227074"
$ws.Range("D20").Value = "RODRIGO ANTONIO VELASQUEZ MORELO"
$ws.Range("E20").Value = "1801"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 737717
$ws.Range("C21").Value = "1050963935"
$ws.Range("D21").Value = "ANGELA INES CAUSIL MARTINEZ"
$ws.Range("E21").Value = "1801"
$ws.Range("F21").Value = 35730
$ws.Range("G21").Value = 893263
$ws.Range("C22").Value = "1050952506"
$ws.Range("D22").Value = "YURIS ROCIO PUELLO OYOLA"
$ws.Range("E22").Value = "1802"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 737717
$ws.Range("C23").Value = "92227074"
$ws.Range("D23").Value = "RODRIGO ANTONIO VELASQUEZ MORELO"
$ws.Range("E23").Value = "1802"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 737717
$ws.Range("C24").Value = "1050963935"
$ws.Range("D24").Value = "ANGELA INES CAUSIL MARTINEZ"
$ws.Range("E24").Value = "1802"
$ws.Range("F24").Value = 35730
$ws.Range("G24").Value = 893263
$ws.Range("C25").Value = "1050952506"
$ws.Range("D25").Value = "YURIS ROCIO PUELLO OYOLA"
$ws.Range("E25").Value = "1803"
$ws.Range("F25").Value = 29509
$ws.Range("G25").Value = 737717
$ws.Range("C26").Value = "92227074"
$ws.Range("D26").Value = "RODRIGO ANTONIO VELASQUEZ MORELO"
$ws.Range("E26").Value = "1803"
$ws.Range("F26").Value = 29509
$ws.Range("G26").Value = 737717
$ws.Range("C27").Value = "1050963935"
$ws.Range("D27").Value = "ANGELA INES CAUSIL MARTINEZ"
$ws.Range("E27").Value = "1803"
$ws.Range("F27").Value = 35730
$ws.Range("G27").Value = 893263
$ws.Range("C28").Value = "1050952506"
$ws.Range("D28").Value = "YURIS ROCIO PUELLO OYOLA"
$ws.Range("E28").Value = "1804"
$ws.Range("F28").Value = 22624
$ws.Range("G28").Value = 737717
$ws.Range("C29").Value = "92227074"
$ws.Range("D29").Value = "RODRIGO ANTONIO VELASQUEZ MORELO"
$ws.Range("E29").Value = "1804"
$ws.Range("F29").Value = 22624
$ws.Range("G29").Value = 737717
$ws.Range("C30").Value = "1050963935"
$ws.Range("D30").Value = "ANGELA INES CAUSIL MARTINEZ"
$ws.Range("E30").Value = "1804"
$ws.Range("F30").Value = 27393
$ws.Range("G30").Value = 893263
